$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename the column headers in row 1
$ws.Range("A1").Value = "randomBalloon"
$ws.Range("B1").Value = "maxPumps"

# Trial data for rows 2-91: balloon color label + max pump value
$data = @(
    @(2, "redBalloon", 4),
    @(3, "greenBalloon", 30),
    @(4, "redBalloon", 3),
    @(5, "greenBalloon", 20),
    @(6, "greenBalloon", 4),
    @(7, "blueBalloon", 39),
    @(8, "greenBalloon", 22),
    @(9, "greenBalloon", 11),
    @(10, "blueBalloon", 36),
    @(11, "blueBalloon", 29),
    @(12, "greenBalloon", 17),
    @(13, "redBalloon", 4),
    @(14, "blueBalloon", 30),
    @(15, "blueBalloon", 21),
    @(16, "blueBalloon", 2),
    @(17, "blueBalloon", 3),
    @(18, "redBalloon", 5),
    @(19, "greenBalloon", 30),
    @(20, "greenBalloon", 22),
    @(21, "redBalloon", 7),
    @(22, "redBalloon", 6),
    @(23, "blueBalloon", 11),
    @(24, "blueBalloon", 20),
    @(25, "greenBalloon", 8),
    @(26, "blueBalloon", 31),
    @(27, "greenBalloon", 6),
    @(28, "redBalloon", 3),
    @(29, "redBalloon", 7),
    @(30, "redBalloon", 4),
    @(31, "redBalloon", 4),
    @(32, "redBalloon", 2),
    @(33, "redBalloon", 2),
    @(34, "redBalloon", 1),
    @(35, "redBalloon", 1),
    @(36, "redBalloon", 7),
    @(37, "redBalloon", 3),
    @(38, "redBalloon", 1),
    @(39, "redBalloon", 7),
    @(40, "redBalloon", 6),
    @(41, "redBalloon", 7),
    @(42, "redBalloon", 1),
    @(43, "redBalloon", 1),
    @(44, "redBalloon", 5),
    @(45, "redBalloon", 1),
    @(46, "redBalloon", 4),
    @(47, "redBalloon", 7),
    @(48, "redBalloon", 5),
    @(49, "redBalloon", 3),
    @(50, "redBalloon", 1),
    @(51, "redBalloon", 2),
    @(52, "greenBalloon", 24),
    @(53, "greenBalloon", 5),
    @(54, "greenBalloon", 4),
    @(55, "greenBalloon", 7),
    @(56, "greenBalloon", 19),
    @(57, "greenBalloon", 25),
    @(58, "greenBalloon", 20),
    @(59, "greenBalloon", 2),
    @(60, "greenBalloon", 14),
    @(61, "greenBalloon", 6),
    @(62, "greenBalloon", 3),
    @(63, "greenBalloon", 20),
    @(64, "greenBalloon", 29),
    @(65, "greenBalloon", 30),
    @(66, "greenBalloon", 3),
    @(67, "greenBalloon", 24),
    @(68, "greenBalloon", 8),
    @(69, "greenBalloon", 4),
    @(70, "greenBalloon", 22),
    @(71, "greenBalloon", 20),
    @(72, "blueBalloon", 26),
    @(73, "blueBalloon", 54),
    @(74, "blueBalloon", 56),
    @(75, "blueBalloon", 54),
    @(76, "blueBalloon", 39),
    @(77, "blueBalloon", 16),
    @(78, "blueBalloon", 56),
    @(79, "blueBalloon", 14),
    @(80, "blueBalloon", 61),
    @(81, "blueBalloon", 9),
    @(82, "blueBalloon", 15),
    @(83, "blueBalloon", 25),
    @(84, "blueBalloon", 52),
    @(85, "blueBalloon", 38),
    @(86, "blueBalloon", 39),
    @(87, "blueBalloon", 55),
    @(88, "blueBalloon", 18),
    @(89, "blueBalloon", 49),
    @(90, "blueBalloon", 53),
    @(91, "blueBalloon", 30)
)

foreach ($row in $data) {
    $r = $row[0]
    $label = $row[1]
    $val = $row[2]
    $ws.Cells.Item($r, 1).Value = $label
    $ws.Cells.Item($r, 2).Value = $val
}
